$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Datos actualizados" timestamp in A1 (03:22 -> 03:52)
$ws.Range("A1").Value = "Datos actualizados a 30 de Abril de 2020 a las 03:52"

# Update Corea del Sur stats (row 38) - plain data refresh, no reordering
$ws.Range("B38").Value = 10765
$ws.Range("C38").Value = 4
$ws.Range("D38").Value = 9059
$ws.Range("E38").Value = 1459
$ws.Range("F38").Value = 55
$ws.Range("G38").Value = 1
$ws.Range("H38").Value = 247

# Panama moves up to row 49 with refreshed stats, Colombia moves down to
# row 50 keeping its previous (unchanged) stats.
$ws.Range("A49").Value = "Panama"
$ws.Range("B49").Value = 6378
$ws.Range("C49").Value = 178
$ws.Range("D49").Value = 527
$ws.Range("E49").Value = 5673
$ws.Range("F49").Value = 92
$ws.Range("G49").Value = 2
$ws.Range("H49").Value = 178

$ws.Range("A50").Value = "Colombia"
$ws.Range("B50").Value = 6207
$ws.Range("C50").Value = 0
$ws.Range("D50").Value = 1411
$ws.Range("E50").Value = 4518
$ws.Range("F50").Value = 118
$ws.Range("G50").Value = 0
$ws.Range("H50").Value = 278
